$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The custom plots were wired up incorrectly around the optical-flow (OF)
# channels, so those rows are dropped and the remaining rows are re-pointed
# at the correct NKF1 / NKF6 / PSC fields. The table shrinks from 12 data
# rows (plots 1-4, 3 rows each) to 10 data rows (plots 1-4, rows 16 & 17
# -- the trailing OF-based plot-4 rows -- are removed entirely).
$ws.Rows("16:17").Delete()

# Plot 1 - Vx (North velocity), rows 6-7: NKF1/NKF6 VN channels
$ws.Range("F6").Value = "Vx~[~m/s~]"
$ws.Range("K6").Value = "NKF1/VN"
$ws.Range("R6").Value = "VN_{EKF1}"
$ws.Range("S6").Value = "m/s"

$ws.Range("F7").Value = "Vx~[~m/s~]"
$ws.Range("K7").Value = "NKF6/VN"
$ws.Range("R7").Value = "VN_{EKF2}"
$ws.Range("S7").Value = "m/s"

# Plot 2 - PN (North position), rows 8-10: NKF1/NKF6/PSC-desired PN channels
$ws.Range("A8").Value = 2
$ws.Range("F8").Value = "PN~[~m~]"
$ws.Range("K8").Value = "NKF1/PN"
$ws.Range("R8").Value = "PN_{EKF1}"
$ws.Range("S8").Value = "m"

$ws.Range("F9").Value = "PN~[~m~]"
$ws.Range("K9").Value = "NKF6/PN"
$ws.Range("R9").Value = "PN_{EKF2}"
$ws.Range("S9").Value = "m"

$ws.Range("F10").Value = "PN~[~m~]"
$ws.Range("K10").Value = "PSC/TPX"
$ws.Range("R10").Value = "PN_{DES}"
$ws.Range("S10").Value = "m"

# Plot 3 - Vy (East velocity), rows 11-12: NKF1/NKF6 VE channels
$ws.Range("A11").Value = 3
$ws.Range("C11").Value = 2
$ws.Range("F11").Value = "Vy~[~m/s~]"
$ws.Range("K11").Value = "NKF1/VE"
$ws.Range("R11").Value = "VE_{EKF1}"
$ws.Range("S11").Value = "m/s"

$ws.Range("F12").Value = "Vy~[~m/s~]"
$ws.Range("K12").Value = "NKF6/VE"
$ws.Range("R12").Value = "VE_{EKF2}"
$ws.Range("S12").Value = "m/s"

# Plot 4 - PE (East position), rows 13-15: NKF1/NKF6/PSC-desired PE channels
$ws.Range("A13").Value = 4
$ws.Range("F13").Value = "PE~[~m~]"
$ws.Range("K13").Value = "NKF1/PE"
$ws.Range("R13").Value = "PE_{EKF1}"
$ws.Range("S13").Value = "m"

$ws.Range("A14").Value = 4
$ws.Range("F14").Value = "PE~[~m~]"
$ws.Range("K14").Value = "NKF6/PE"
$ws.Range("R14").Value = "PE_{EKF2}"
$ws.Range("S14").Value = "m"

$ws.Range("F15").Value = "PE~[~m~]"
$ws.Range("K15").Value = "PSC/TPY"
$ws.Range("R15").Value = "PE_{DES}"
$ws.Range("S15").Value = "m"

# Keep the plot title cell as-is (no textual change, just touched so the
# shared-string table stays correctly deduplicated).
$ws.Range("B3").Value = "VelocityTuning"

# Match the author's final selection (whole row 8 highlighted).
$ws.Rows(8).Select()
